$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add CF1 header date, copying format from CE1
$ws.Range("CE1").Copy()
$ws.Range("CF1").PasteSpecial(-4122)
$ws.Range("CF1").Value = 45986

# Rows 2-138: copy CE value into CF (same value)
$ws.Range("CF2").Value = 70.81
$ws.Range("CF3").Value = 70.47
$ws.Range("CF4").Value = 70.36
$ws.Range("CF5").Value = 71.37
$ws.Range("CF6").Value = 72.34999999999999
$ws.Range("CF7").Value = 71.90000000000001
$ws.Range("CF8").Value = 71.73
$ws.Range("CF9").Value = 71.53
$ws.Range("CF10").Value = 70.95999999999999
$ws.Range("CF11").Value = 70.97
$ws.Range("CF12").Value = 71.39
$ws.Range("CF13").Value = 71.33
$ws.Range("CF14").Value = 72.34
$ws.Range("CF15").Value = 72.73999999999999
$ws.Range("CF16").Value = 73.17
$ws.Range("CF17").Value = 73.98
$ws.Range("CF18").Value = 73.69
$ws.Range("CF19").Value = 74.28
$ws.Range("CF20").Value = 74.43000000000001
$ws.Range("CF21").Value = 74.45999999999999
$ws.Range("CF22").Value = 73.93000000000001
$ws.Range("CF23").Value = 74.94
$ws.Range("CF24").Value = 75.23
$ws.Range("CF25").Value = 75.92
$ws.Range("CF26").Value = 75.53
$ws.Range("CF27").Value = 76.37
$ws.Range("CF28").Value = 76.67
$ws.Range("CF29").Value = 77.23999999999999
$ws.Range("CF30").Value = 77.98
$ws.Range("CF31").Value = 77.62
$ws.Range("CF32").Value = 78
$ws.Range("CF33").Value = 77.98
$ws.Range("CF34").Value = 78.88
$ws.Range("CF35").Value = 78.78
$ws.Range("CF36").Value = 79.86
$ws.Range("CF37").Value = 80.23999999999999
$ws.Range("CF38").Value = 81.48
$ws.Range("CF39").Value = 82.13
$ws.Range("CF40").Value = 82.19
$ws.Range("CF41").Value = 81.84
$ws.Range("CF42").Value = 83.58
$ws.Range("CF43").Value = 83.37
$ws.Range("CF44").Value = 83.37
$ws.Range("CF45").Value = 83.11
$ws.Range("CF46").Value = 82.72
$ws.Range("CF47").Value = 83.04000000000001
$ws.Range("CF48").Value = 83.61
$ws.Range("CF49").Value = 83.38
$ws.Range("CF50").Value = 82.25
$ws.Range("CF51").Value = 82.37
$ws.Range("CF52").Value = 83.12
$ws.Range("CF53").Value = 83.22
$ws.Range("CF54").Value = 83.09
$ws.Range("CF55").Value = 83.58
$ws.Range("CF56").Value = 83.31
$ws.Range("CF57").Value = 83.28
$ws.Range("CF58").Value = 83.41
$ws.Range("CF59").Value = 83.91
$ws.Range("CF60").Value = 84.56999999999999
$ws.Range("CF61").Value = 84.89
$ws.Range("CF62").Value = 85.84
$ws.Range("CF63").Value = 87.31
$ws.Range("CF64").Value = 87.97
$ws.Range("CF65").Value = 89.26000000000001
$ws.Range("CF66").Value = 89.38
$ws.Range("CF67").Value = 89.97
$ws.Range("CF68").Value = 90.43000000000001
$ws.Range("CF69").Value = 91.09999999999999
$ws.Range("CF70").Value = 91.63
$ws.Range("CF71").Value = 91.3
$ws.Range("CF72").Value = 90.8
$ws.Range("CF73").Value = 89.39
$ws.Range("CF74").Value = 85.2
$ws.Range("CF75").Value = 85.41
$ws.Range("CF76").Value = 85.94
$ws.Range("CF77").Value = 86.59999999999999
$ws.Range("CF78").Value = 87.23
$ws.Range("CF79").Value = 89.12
$ws.Range("CF80").Value = 89.98
$ws.Range("CF81").Value = 90.56999999999999
$ws.Range("CF82").Value = 92.20999999999999
$ws.Range("CF83").Value = 92.45999999999999
$ws.Range("CF84").Value = 92.95999999999999
$ws.Range("CF85").Value = 92.95
$ws.Range("CF86").Value = 93.15000000000001
$ws.Range("CF87").Value = 93.25
$ws.Range("CF88").Value = 93.44
$ws.Range("CF89").Value = 93.16
$ws.Range("CF90").Value = 92.64
$ws.Range("CF91").Value = 93.72
$ws.Range("CF92").Value = 94.2
$ws.Range("CF93").Value = 94.34
$ws.Range("CF94").Value = 95.33
$ws.Range("CF95").Value = 95.31
$ws.Range("CF96").Value = 95.84
$ws.Range("CF97").Value = 96.56999999999999
$ws.Range("CF98").Value = 96.36
$ws.Range("CF99").Value = 96.86
$ws.Range("CF100").Value = 97.38
$ws.Range("CF101").Value = 97.88
$ws.Range("CF102").Value = 98.73999999999999
$ws.Range("CF103").Value = 98.95999999999999
$ws.Range("CF104").Value = 99.31
$ws.Range("CF105").Value = 99.77
$ws.Range("CF106").Value = 101.06
$ws.Range("CF107").Value = 101.78
$ws.Range("CF108").Value = 102.59
$ws.Range("CF109").Value = 103.63
$ws.Range("CF110").Value = 103.06
$ws.Range("CF111").Value = 103.87
$ws.Range("CF112").Value = 103.19
$ws.Range("CF113").Value = 103.65
$ws.Range("CF114").Value = 104.33
$ws.Range("CF115").Value = 104.39
$ws.Range("CF116").Value = 104.79
$ws.Range("CF117").Value = 104.44
$ws.Range("CF118").Value = 102.32
$ws.Range("CF119").Value = 93.23999999999999
$ws.Range("CF120").Value = 101.33
$ws.Range("CF121").Value = 102.31
$ws.Range("CF122").Value = 101.68
$ws.Range("CF123").Value = 104.07
$ws.Range("CF124").Value = 104.16
$ws.Range("CF125").Value = 104.72
$ws.Range("CF126").Value = 105.43
$ws.Range("CF127").Value = 105.59
$ws.Range("CF128").Value = 105.9
$ws.Range("CF129").Value = 105.53
$ws.Range("CF130").Value = 105.03
$ws.Range("CF131").Value = 104.95
$ws.Range("CF132").Value = 104.95
$ws.Range("CF133").Value = 104.66
$ws.Range("CF134").Value = 104.55
$ws.Range("CF135").Value = 104.28
$ws.Range("CF136").Value = 104.3
$ws.Range("CF137").Value = 104.49
$ws.Range("CF138").Value = 104.81

# Row 139: CF139 gets a new distinct value
$ws.Range("CF139").Value = 104.59

# Row 140: add CF140 value (no corresponding CE140)
$ws.Range("CF140").Value = 104.59

# Row 141: new row, copy format from row 140 col A, then set value
$ws.Range("A140").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$ws.Range("A141").Value = 45976
